$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 5726.7896
$ws.Range("I33").Value = 57.416668
$ws.Range("J33").Value = 15445.714
$ws.Range("K33").Value = 57.416668
$ws.Range("L33").Value = 15445.714
$ws.Range("M33").Value = 171.583332
$ws.Range("N33").Value = -15903.714
$ws.Range("H137").Value = 804.1591
$ws.Range("I137").Value = 712.9583
$ws.Range("J137").Value = 913.6
$ws.Range("K137").Value = 2138.8749
$ws.Range("L137").Value = 2740.8
$ws.Range("M137").Value = 411.1251000000002
$ws.Range("N137").Value = -7840.8
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1445.9642
$ws.Range("I61").Value = 911.0714
$ws.Range("J61").Value = 1980.8572
$ws.Range("K61").Value = 911.0714
$ws.Range("L61").Value = 1980.8572
$ws.Range("M61").Value = -699.0714
$ws.Range("N61").Value = -2404.8572
$ws.Range("H74").Value = 1235.9333
$ws.Range("I74").Value = 1095.238
$ws.Range("J74").Value = 1564.2222
$ws.Range("K74").Value = 1095.238
$ws.Range("L74").Value = 1564.2222
$ws.Range("M74").Value = -221.2380000000001
$ws.Range("N74").Value = -3312.2222
$ws.Range("H77").Value = 1235.9333
$ws.Range("I77").Value = 1095.238
$ws.Range("J77").Value = 1564.2222
$ws.Range("K77").Value = 5476.190000000001
$ws.Range("L77").Value = 7821.111
$ws.Range("M77").Value = -1108.190000000001
$ws.Range("N77").Value = -16557.111
$ws.Range("H132").Value = 1127.5
$ws.Range("I132").Value = 646.13513
$ws.Range("J132").Value = 3106.4443
$ws.Range("K132").Value = 1938.40539
$ws.Range("L132").Value = 9319.332900000001
$ws.Range("M132").Value = 591.5946100000001
$ws.Range("N132").Value = -14379.3329
$ws.Range("H136").Value = 1445.9642
$ws.Range("I136").Value = 911.0714
$ws.Range("J136").Value = 1980.8572
$ws.Range("K136").Value = 2733.2142
$ws.Range("L136").Value = 5942.571599999999
$ws.Range("M136").Value = -183.2142000000003
$ws.Range("N136").Value = -11042.5716
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 12210.182
$ws.Range("I75").Value = 10590.333
$ws.Range("J75").Value = 19499.5
$ws.Range("K75").Value = 10590.333
$ws.Range("L75").Value = 19499.5
$ws.Range("M75").Value = -9654.333000000001
$ws.Range("N75").Value = -21371.5
$ws.Range("H78").Value = 12210.182
$ws.Range("I78").Value = 10590.333
$ws.Range("J78").Value = 19499.5
$ws.Range("K78").Value = 31770.999
$ws.Range("L78").Value = 58498.5
$ws.Range("M78").Value = -27090.999
$ws.Range("N78").Value = -67858.5
$ws.Range("H80").Value = 1956.6316
$ws.Range("I80").Value = 277.33334
$ws.Range("K80").Value = 277.33334
$ws.Range("M80").Value = 720.66666
$ws.Range("H83").Value = 1956.6316
$ws.Range("I83").Value = 277.33334
$ws.Range("K83").Value = 1386.6667
$ws.Range("M83").Value = 3605.3333
$ws.Range("H86").Value = 34487530
$ws.Range("I86").Value = 55558172
$ws.Range("J86").Value = 8287.909
$ws.Range("K86").Value = 55558172
$ws.Range("L86").Value = 8287.909
$ws.Range("M86").Value = -55557049
$ws.Range("N86").Value = -10533.909
$ws.Range("H89").Value = 34487530
$ws.Range("I89").Value = 55558172
$ws.Range("J89").Value = 8287.909
$ws.Range("K89").Value = 277790860
$ws.Range("L89").Value = 41439.545
$ws.Range("M89").Value = -277785244
$ws.Range("N89").Value = -52671.545
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H64").Value = 9966.666999999999
$ws.Range("J64").Value = 9966.666999999999
$ws.Range("L64").Value = 9966.666999999999
$ws.Range("N64").Value = -10462.667
$ws.Range("H67").Value = 9966.666999999999
$ws.Range("J67").Value = 9966.666999999999
$ws.Range("L67").Value = 9966.666999999999
$ws.Range("N67").Value = -11682.667
$ws.Range("H119").Value = 50000
$ws.Range("J119").Value = 50000
$ws.Range("L119").Value = 50000
$ws.Range("N119").Value = -59676
$ws.Range("H132").Value = 54108.684
$ws.Range("I132").Value = 1210.7858
$ws.Range("K132").Value = 3632.3574
$ws.Range("M132").Value = -1102.3574
$ws.Range("H134").Value = 23329.4
$ws.Range("I134").Value = 24912.453
$ws.Range("J134").Value = 1166.6666
$ws.Range("K134").Value = 74737.359
$ws.Range("L134").Value = 3499.9998
$ws.Range("M134").Value = -72202.359
$ws.Range("N134").Value = -8569.9998
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 140.95454
$ws.Range("I33").Value = 50.29032
$ws.Range("J33").Value = 357.15384
$ws.Range("K33").Value = 301.74192
$ws.Range("L33").Value = 2142.92304
$ws.Range("M33").Value = -18.74191999999999
$ws.Range("N33").Value = -2708.92304
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 169252.33
$ws.Range("I44").Value = 964
$ws.Range("J44").Value = 202910
$ws.Range("K44").Value = 964
$ws.Range("L44").Value = 202910
$ws.Range("M44").Value = -368
$ws.Range("N44").Value = -204102
$ws.Range("H102").Value = 1711.7826
$ws.Range("I102").Value = 1675.05
$ws.Range("J102").Value = 1956.6666
$ws.Range("K102").Value = 1675.05
$ws.Range("L102").Value = 1956.6666
$ws.Range("M102").Value = -53.04999999999995
$ws.Range("N102").Value = -5200.6666
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 606.5469000000001
$ws.Range("I22").Value = 513.9535
$ws.Range("J22").Value = 796.1429000000001
$ws.Range("K22").Value = 513.9535
$ws.Range("L22").Value = 796.1429000000001
$ws.Range("M22").Value = -218.9535
$ws.Range("N22").Value = -1386.1429
$ws.Range("H27").Value = 606.5469000000001
$ws.Range("I27").Value = 513.9535
$ws.Range("J27").Value = 796.1429000000001
$ws.Range("K27").Value = 513.9535
$ws.Range("L27").Value = 796.1429000000001
$ws.Range("M27").Value = -406.9535
$ws.Range("N27").Value = -1010.1429
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
$ws.Range("H132").Value = 1672.0112
$ws.Range("I132").Value = 1851.2153
$ws.Range("J132").Value = 1186.6666
$ws.Range("K132").Value = 5553.6459
$ws.Range("L132").Value = 3559.9998
$ws.Range("M132").Value = -3023.6459
$ws.Range("N132").Value = -8619.9998
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 87400
$ws.Range("I96").Value = 1975
$ws.Range("J96").Value = 130112.5
$ws.Range("K96").Value = 1975
$ws.Range("L96").Value = 130112.5
$ws.Range("M96").Value = -602
$ws.Range("N96").Value = -132858.5
$ws.Range("H121").Value = 31000
$ws.Range("J121").Value = 31000
$ws.Range("L121").Value = 31000
$ws.Range("N121").Value = -34494
$ws.Range("H122").Value = 1984.3684
$ws.Range("I122").Value = 1743.0714
$ws.Range("J122").Value = 2660
$ws.Range("K122").Value = 5229.2142
$ws.Range("L122").Value = 7980
$ws.Range("M122").Value = -2779.2142
$ws.Range("N122").Value = -12880
$ws.Range("H123").Value = 49980
$ws.Range("J123").Value = 49980
$ws.Range("L123").Value = 49980
$ws.Range("N123").Value = -59780
